# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).
#
# Changes:
#   F2: 181  -> 182
#   F7: 5738 -> 5742
#   F8: 20   -> 22
# applied identically on sheet "展览" (row/col based data) and sheet "全部类型".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 182
    $ws.Range("F7").Value = 5742
    $ws.Range("F8").Value = 22
}
